$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at 30 (shifts old rows 30-44 down to 31-45,
#    and auto-adjusts formula references to rows >= 30).
$ws.Rows.Item(30).Insert()

# 2. Populate the newly inserted row 30 ("seqAttackPenalty") with static zero values.
$ws.Range("A30").Value = "seqAttackPenalty"
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0

# 3. Update row 31 ("score") formula to include the new (1+B30) factor.
$ws.Range("B31").Formula = "=MAX(-100, IF(B28>B`$26, POWER(B29, 3 + B26), -POWER(B29, 4+B26)*(1+B30)))"
$ws.Range("C31").Formula = "=MAX(-100, IF(C28>C`$26, POWER(C29, 3 + C26), -POWER(C29, 4+C26)*(1+C30)))"
$ws.Range("D31").Formula = "=MAX(-100, IF(D28>D`$26, POWER(D29, 3 + D26), -POWER(D29, 4+D26)*(1+D30)))"

# 4. Row 32 previously held a duplicate/buggy "score" formula (referencing nonexistent row 27).
#    Replace it with the same corrected formula used in row 31.
$ws.Range("B32").Formula = "=MAX(-100, IF(B28>B`$26, POWER(B29, 3 + B26), -POWER(B29, 4+B26)*(1+B30)))"
$ws.Range("C32").Formula = "=MAX(-100, IF(C28>C`$26, POWER(C29, 3 + C26), -POWER(C29, 4+C26)*(1+C30)))"
$ws.Range("D32").Formula = "=MAX(-100, IF(D28>D`$26, POWER(D29, 3 + D26), -POWER(D29, 4+D26)*(1+D30)))"

# Row 32 lacked the aggregate columns E/F/H/I/J/K that row 31 has; add matching ones.
$ws.Range("E32").Value = "total:"
$ws.Range("F32").Formula = "=MAX(-1000, SUM(B32:D32))"
$ws.Range("H32").Value = "certainty =>"
$ws.Range("I32").Formula = "= MAX(0, AVERAGE(B28:D28) - _xlfn.STDEV.P(B28:D28) * 2)"
$ws.Range("J32").Value = "weight =>"
$ws.Range("K32").Formula = "=MAX((1/(1+EXP(-F32/100)))*10-4, 0) * I32"

# 5. Update the two non-i.i.d. data input rows.
$ws.Range("B26").Value = 1
$ws.Range("B28").Value = 1
$ws.Range("D28").Value = 0.8

# 6. Update the active view/selection to match what was saved.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select()

$wb.Application.Calculate()
